# Latest data from Bultark
# Update the year-to-date poker figures (Sheet1, rows 191-199) with the
# latest month's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 191 - Richard
$ws.Cells.Item(191, 4).Value = 41
$ws.Cells.Item(191, 6).Value = 41
$ws.Cells.Item(191, 7).Value = 127950
$ws.Cells.Item(191, 8).Value = 150

# Row 192 - Mark
$ws.Cells.Item(192, 4).Value = 34
$ws.Cells.Item(192, 6).Value = 34
$ws.Cells.Item(192, 7).Value = 102700
$ws.Cells.Item(192, 8).Value = 130
$ws.Cells.Item(192, 9).Value = 60

# Row 193 - Andy
$ws.Cells.Item(193, 4).Value = 26
$ws.Cells.Item(193, 6).Value = 26
$ws.Cells.Item(193, 7).Value = 82450
$ws.Cells.Item(193, 9).Value = 0

# Row 194 - Anthony (unchanged)

# Row 195 - now Pepe (was Jon)
$ws.Cells.Item(195, 2).Value = "Pepe"
$ws.Cells.Item(195, 4).Value = 22
$ws.Cells.Item(195, 6).Value = 22
$ws.Cells.Item(195, 7).Value = 62350
$ws.Cells.Item(195, 8).Value = 40
$ws.Cells.Item(195, 9).Value = -20
$ws.Cells.Item(195, 11).Value = 364

# Row 196 - now Matt (was Pepe)
$ws.Cells.Item(196, 2).Value = "Matt"
$ws.Cells.Item(196, 4).Value = 19
$ws.Cells.Item(196, 6).Value = 19
$ws.Cells.Item(196, 7).Value = 77150
$ws.Cells.Item(196, 8).Value = 30
$ws.Cells.Item(196, 9).Value = -40
$ws.Cells.Item(196, 11).Value = 362

# Row 197 - now Jon (was Prashant)
$ws.Cells.Item(197, 2).Value = "Jon"
$ws.Cells.Item(197, 4).Value = 19
$ws.Cells.Item(197, 6).Value = 19
$ws.Cells.Item(197, 7).Value = 56750
$ws.Cells.Item(197, 8).Value = 20
$ws.Cells.Item(197, 9).Value = -50
$ws.Cells.Item(197, 11).Value = 357

# Row 198 - now Prashant (was Matt)
$ws.Cells.Item(198, 2).Value = "Prashant"
$ws.Cells.Item(198, 4).Value = 14
$ws.Cells.Item(198, 6).Value = 14
$ws.Cells.Item(198, 7).Value = 49350
$ws.Cells.Item(198, 8).Value = 60
$ws.Cells.Item(198, 9).Value = 10
$ws.Cells.Item(198, 11).Value = 365

# Row 199 - Maisy (only Takehome changes)
$ws.Cells.Item(199, 9).Value = -40

$wb.Save()
